$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '61.890.78'
$ws.Range('E2').Value = '  -2.01%  '

# Row 3
$ws.Range('D3').Value = '3.416.97'
$ws.Range('E3').Value = '  -1.80%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
Set-TextValue $ws.Range('D5') '406.24'
$ws.Range('E5').Value = '  -0.89%  '

# Row 6
Set-TextValue $ws.Range('D6') '132.73'
$ws.Range('E6').Value = '  -0.25%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.592'
$ws.Range('E7').Value = '  -2.18%  '

# Row 8
$ws.Range('E8').Value = '  +0.05%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.685'
$ws.Range('E9').Value = '  -1.13%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.126'
$ws.Range('E10').Value = '  -3.29%  '

# Row 11
Set-TextValue $ws.Range('D11') '41.94'
$ws.Range('E11').Value = '  -2.90%  '

# Row 12
$ws.Range('E12').Value = '  -1.03%  '

# Row 13
Set-TextValue $ws.Range('D13') '8.44'
$ws.Range('E13').Value = '  -4.01%  '

# Row 14
Set-TextValue $ws.Range('D14') '19.79'
$ws.Range('E14').Value = '  -2.19%  '

# Row 15
$ws.Range('D15').Value = '3.425.40'
$ws.Range('E15').Value = '  -2.04%  '

# Row 16
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D16') '11.56'
$ws.Range('E16').Value = '  +6.38%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '61.983.10'
$ws.Range('E17').Value = '  -2.02%  '

# Row 18
$ws.Range('E18').Value = '  -2.94%  '

# Row 19
Set-TextValue $ws.Range('D19') '0.0000143'
$ws.Range('E19').Value = '  +1.87%  '

# Row 20
Set-TextValue $ws.Range('D20') '3.16'
$ws.Range('E20').Value = '  -5.46%  '

# Row 21
Set-TextValue $ws.Range('D21') '83.48'
$ws.Range('E21').Value = '  +0.96%  '

# Row 22
Set-TextValue $ws.Range('D22') '310.87'
$ws.Range('E22').Value = '  -0.86%  '

# Row 23
Set-TextValue $ws.Range('D23') '12.82'
$ws.Range('E23').Value = '  -2.80%  '

# Row 24
Set-TextValue $ws.Range('D24') '3.14'
$ws.Range('E24').Value = '  -1.21%  '

# Row 25
Set-TextValue $ws.Range('D25') '4.80'
$ws.Range('E25').Value = '  +9.10%  '

# Row 26
Set-TextValue $ws.Range('D26') '29.64'
$ws.Range('E26').Value = '  -2.92%  '

# Row 27
Set-TextValue $ws.Range('D27') '8.11'
$ws.Range('E27').Value = '  -1.24%  '

# Row 28
Set-TextValue $ws.Range('D28') '7.69'
$ws.Range('E28').Value = '  +0.39%  '

# Row 29
Set-TextValue $ws.Range('D29') '2.76'
$ws.Range('E29').Value = '  +4.64%  '

# Row 30
$ws.Range('E30').Value = '  -2.69%  '

# Row 31
$ws.Range('E31').Value = '  -3.51%  '

# Row 32
$ws.Range('E32').Value = '  -3.50%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.998'
$ws.Range('E33').Value = '  -0.06%  '

# Row 34
Set-TextValue $ws.Range('D34') '11.35'
$ws.Range('E34').Value = '  -3.83%  '

# Row 35
Set-TextValue $ws.Range('D35') '0.0484'
$ws.Range('E35').Value = '  -1.84%  '

# Row 36
Set-TextValue $ws.Range('D36') '51.14'
$ws.Range('E36').Value = '  -2.79%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.999'
$ws.Range('E37').Value = '  -0.05%  '

# Row 38
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D38') '0.327'
$ws.Range('E38').Value = '  +13.43%  '

# Row 39
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D39') '3.38'
$ws.Range('E39').Value = '  -5.81%  '

# Row 40
$ws.Range('E40').Value = '  -3.91%  '

# Row 41
Set-TextValue $ws.Range('D41') '138.30'
$ws.Range('E41').Value = '  +1.54%  '

# Row 42
$ws.Range('E42').Value = '  -0.57%  '

# Row 43
$ws.Range('E43').Value = '  -1.05%  '

# Row 44
Set-TextValue $ws.Range('D44') '3.96'
$ws.Range('E44').Value = '  -0.73%  '

# Row 45
Set-TextValue $ws.Range('D45') '16.70'
$ws.Range('E45').Value = '  -4.06%  '

# Row 46
Set-TextValue $ws.Range('D46') '2.22'
$ws.Range('E46').Value = '  -1.21%  '

# Row 47
Set-TextValue $ws.Range('D47') '21.21'
$ws.Range('E47').Value = '  -5.02%  '

# Row 48
$ws.Range('D48').Value = '2.117.24'
$ws.Range('E48').Value = '  -3.68%  '

# Row 49
$ws.Range('E49').Value = '  -3.88%  '

# Row 50
Set-TextValue $ws.Range('D50') '1.77'
$ws.Range('E50').Value = '  +21.08%  '

# Row 51
Set-TextValue $ws.Range('D51') '1.93'
$ws.Range('E51').Value = '  +2.60%  '
